$wb = $excel.ActiveWorkbook

$updates = @{
    "weibull" = @{
        "B2" = -3.92419203737332
        "C2" = 0.494772230509085
        "B3" = 0.360426783558681
        "C3" = 0.187127507302931
    }
    "lognormal" = @{
        "B2" = 2.81513530583061
        "C2" = 0.549529613637574
        "B3" = -1.03564365282378
        "C3" = 0.177278779745409
    }
    "llogis" = @{
        "B2" = -2.76933170304976
        "C2" = 0.177956593636583
        "B3" = 1.91020861268176
        "C3" = 0.371607190811063
    }
    "gompertz" = @{
        "B2" = -3.42672233937081
        "C2" = 0.243199656921414
        "B3" = 0.0243636821714332
        "C3" = 0.0122956675768103
    }
    "weibull cov" = @{
        "A2" = 0.244799560082935
        "B2" = -0.0810677652745029
        "A3" = -0.0810677652745029
        "B3" = 0.0350167039894086
    }
    "lognormal cov" = @{
        "A2" = 0.301982796264662
        "B2" = -0.0905067219448786
        "A3" = -0.0905067219448786
        "B3" = 0.0314277657480212
    }
    "llogis cov" = @{
        "A2" = 0.0316685492187359
        "B2" = -0.0048304199862856
        "A3" = -0.0048304199862856
        "B3" = 0.13809190426249
    }
    "gompertz cov" = @{
        "A2" = 0.0591460731266934
        "B2" = -0.00126953734763064
        "A3" = -0.00126953734763064
        "B3" = 0.000151183441159423
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
